$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add column A values
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5

# Update column B base value and formulas (shift by 152000 seconds)
$ws.Range("B3").Value = 1592656000
$ws.Range("B4").Formula = "=B3+3600"
$ws.Range("B5").Formula = "=B4+3600"
$ws.Range("B6").Formula = "=B5+3600"
$ws.Range("B7").Formula = "=B6+3600"
$ws.Range("B8").Formula = "=B7+3600"
